# Automatische test-sync: 2025-08-03 15:08:50
# Adds the new Testmail #17 log entry (row 25) to the "Logs" sheet,
# extends the conditional formatting ranges to include it, and bumps
# the "Planning / Afspraak" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 25

$ws.Cells.Item($newRow, 1).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value = "Testmail #17: Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item($newRow, 4).Value = "Planning / Afspraak"
$ws.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Cells.Item($newRow, 6).Value = "2025-08-03 15:08:38"
$ws.Cells.Item($newRow, 7).Value = "Ja"
$ws.Cells.Item($newRow, 8).Value = "Ja"
$ws.Cells.Item($newRow, 9).Value = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional-formatting sqref ranges (D/G/H/I/J 2:24 -> 2:25)
# so the newly added row is covered too. Modifying one rule in a grouped
# conditional-formatting block updates the whole shared sqref.
$ws.Range("D2:D24").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D25"))
$ws.Range("G2:G24").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G25"))
$ws.Range("H2:H24").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H25"))
$ws.Range("I2:I24").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I25"))
$ws.Range("J2:J24").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J25"))

# Update the Dashboard summary count for "Planning / Afspraak" (5 -> 6)
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 2).Value = 6
